# Auto-generated edit script for Maduin_Profits workbook update
# Applies per-cell numeric updates to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1750
$ws.Range("J17").Value = 1750
$ws.Range("L17").Value = 5250
$ws.Range("N17").Value = -5586
$ws.Range("H32").Value = 1169.1666
$ws.Range("I32").Value = 1233.6666
$ws.Range("J32").Value = 1104.6666
$ws.Range("K32").Value = 1233.6666
$ws.Range("L32").Value = 1104.6666
$ws.Range("M32").Value = -907.6666
$ws.Range("N32").Value = -1756.6666
$ws.Range("H53").Value = 529.8
$ws.Range("J53").Value = 875
$ws.Range("L53").Value = 875
$ws.Range("N53").Value = -2149
$ws.Range("H103").Value = 2999
$ws.Range("J103").Value = 3357.8
$ws.Range("L103").Value = 10073.4
$ws.Range("N103").Value = -11245.4
$ws.Range("H112").Value = 3108.375
$ws.Range("J112").Value = 3222.4
$ws.Range("L112").Value = 9667.200000000001
$ws.Range("N112").Value = -11883.2
$ws.Range("H138").Value = 3801.9565
$ws.Range("J138").Value = 3731.3684
$ws.Range("L138").Value = 11194.1052
$ws.Range("N138").Value = -21474.1052

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 8250
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 8250
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 8250
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -8590
$ws.Range("H132").Value = 2056.1428
$ws.Range("I132").Value = 2198
$ws.Range("J132").Value = 1999.4
$ws.Range("K132").Value = 6594
$ws.Range("L132").Value = 5998.200000000001
$ws.Range("M132").Value = -4064
$ws.Range("N132").Value = -11058.2

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 595.7778
$ws.Range("I107").Value = 519.8333
$ws.Range("J107").Value = 747.6667
$ws.Range("K107").Value = 519.8333
$ws.Range("L107").Value = 747.6667
$ws.Range("M107").Value = 1400.1667
$ws.Range("N107").Value = -4587.6667
$ws.Range("H119").Value = 9585.333000000001
$ws.Range("J119").Value = 9585.333000000001
$ws.Range("L119").Value = 9585.333000000001
$ws.Range("N119").Value = -19261.333

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 357.22223
$ws.Range("I5").Value = 172.42857
$ws.Range("K5").Value = 172.42857
$ws.Range("M5").Value = -60.42857000000001
$ws.Range("H6").Value = 5398
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 6497.5
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 6497.5
$ws.Range("M6").Value = -887
$ws.Range("N6").Value = -6723.5
$ws.Range("H10").Value = 1150
$ws.Range("I10").Value = 1150
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1150
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -1011
$ws.Range("N10").ClearContents()
$ws.Range("H59").Value = 37524.75
$ws.Range("I59").Value = 25051
$ws.Range("J59").Value = 49998.5
$ws.Range("K59").Value = 25051
$ws.Range("L59").Value = 49998.5
$ws.Range("M59").Value = -23906
$ws.Range("N59").Value = -52288.5
$ws.Range("H60").Value = 20359
$ws.Range("I60").Value = 9000
$ws.Range("J60").Value = 27931.666
$ws.Range("K60").Value = 9000
$ws.Range("L60").Value = 27931.666
$ws.Range("M60").Value = -8489
$ws.Range("N60").Value = -28953.666

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H107").Value = 935.2222
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 1123.4
$ws.Range("K107").Value = 2100
$ws.Range("L107").Value = 3370.2
$ws.Range("M107").Value = -180
$ws.Range("N107").Value = -7210.200000000001

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2500833.2
$ws.Range("J3").Value = 2250
$ws.Range("L3").Value = 2250
$ws.Range("N3").Value = -2482
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H43").Value = 42499.5
$ws.Range("J43").Value = 42499.5
$ws.Range("L43").Value = 42499.5
$ws.Range("N43").Value = -42801.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 2000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 2000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -2226
$ws.Range("H22").Value = 2006.9166
$ws.Range("I22").Value = 1197.875
$ws.Range("J22").Value = 3625
$ws.Range("K22").Value = 1197.875
$ws.Range("L22").Value = 3625
$ws.Range("M22").Value = -902.875
$ws.Range("N22").Value = -4215
$ws.Range("H26").Value = 4009
$ws.Range("I26").Value = 4009
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 4009
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -3714
$ws.Range("N26").ClearContents()
$ws.Range("H27").Value = 2006.9166
$ws.Range("I27").Value = 1197.875
$ws.Range("J27").Value = 3625
$ws.Range("K27").Value = 1197.875
$ws.Range("L27").Value = 3625
$ws.Range("M27").Value = -1090.875
$ws.Range("N27").Value = -3839
$ws.Range("H28").Value = 2000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 2000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 2000
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -2464
$ws.Range("H29").Value = 27499.5
$ws.Range("J29").Value = 29999
$ws.Range("L29").Value = 29999
$ws.Range("N29").Value = -30589
$ws.Range("H37").Value = 2000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 2000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 2000
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -2214
$ws.Range("H82").Value = 3749.9
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 3749.9
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1745
$ws.Range("J8").Value = 1745
$ws.Range("L8").Value = 1745
$ws.Range("N8").Value = -2025
$ws.Range("H13").Value = 4569.4287
$ws.Range("J13").Value = 4569.4287
$ws.Range("L13").Value = 4569.4287
$ws.Range("N13").Value = -4849.4287
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H17").Value = 1545
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 2990
$ws.Range("K17").Value = 100
$ws.Range("L17").Value = 2990
$ws.Range("M17").Value = 72
$ws.Range("N17").Value = -3334
$ws.Range("H113").Value = 1433
$ws.Range("I113").Value = 1232
$ws.Range("K113").Value = 3696
$ws.Range("M113").Value = -1526
$ws.Range("H125").Value = 59950
$ws.Range("J125").Value = 59950
$ws.Range("L125").Value = 59950
$ws.Range("N125").Value = -69790
